$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06762866666666667
$ws.Range("H2").Value = 0.202886
$ws.Range("I2").Value = 0.0134153952845566
$ws.Range("J2").Value = 0.0134153952845566
$ws.Range("M2").Value = 1.824475333333333
$ws.Range("N2").Value = 5.473426
$ws.Range("O2").Value = 0.1906606574278047
$ws.Range("P2").Value = 0.2015451970524477
$ws.Range("Q2").Value = 0.1233868341595556
$ws.Range("R2").Value = 1.110481507436
$ws.Range("S2").Value = 0.002557788084607433
$ws.Range("T2").Value = 0.002703808486162437
$ws.Range("G3").Value = 0.06762866666666667
$ws.Range("H3").Value = 0.202886
$ws.Range("I3").Value = 0.0134153952845566
$ws.Range("J3").Value = 0.0134153952845566
$ws.Range("O3").Value = 0.6423822165107047
$ws.Range("P3").Value = 0.6790548829333741
$ws.Range("Q3").Value = 0.4157203121240001
$ws.Range("R3").Value = 3.741482809116
$ws.Range("S3").Value = 0.008617811358260724
$ws.Range("T3").Value = 0.009109789674459521
$ws.Range("G4").Value = 0.06762866666666667
$ws.Range("H4").Value = 0.202886
$ws.Range("I4").Value = 0.0134153952845566
$ws.Range("J4").Value = 0.0134153952845566
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009389666666666666
$ws.Range("N4").Value = 0.028169
$ws.Range("O4").Value = 0.0009812355294625031
$ws.Range("P4").Value = 0.001037252838673693
$ws.Range("Q4").Value = 0.0006350106371111111
$ws.Range("R4").Value = 0.005715095734
$ws.Range("S4").Value = 0.00001316366249499066
$ws.Range("T4").Value = 0.00001391515684083601
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.06762866666666667
$ws.Range("H5").Value = 0.202886
$ws.Range("I5").Value = 0.0134153952845566
$ws.Range("J5").Value = 0.0134153952845566
$ws.Range("M5").Value = 1.5503715
$ws.Range("N5").Value = 3.100743
$ws.Range("O5").Value = 0.1620163583726162
$ws.Range("P5").Value = 0.1141770910840848
$ws.Range("Q5").Value = 0.104849557383
$ws.Range("R5").Value = 0.629097344298
$ws.Range("S5").Value = 0.002173513490133028
$ws.Range("T5").Value = 0.00153173080933382
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 0.06762866666666667
$ws.Range("H6").Value = 0.202886
$ws.Range("I6").Value = 0.0134153952845566
$ws.Range("J6").Value = 0.0134153952845566
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03788966666666667
$ws.Range("N6").Value = 0.113669
$ws.Range("O6").Value = 0.003959532159411881
$ws.Range("P6").Value = 0.004185576091419648
$ws.Range("Q6").Value = 0.002562427637111111
$ws.Range("R6").Value = 0.023061848734
$ws.Range("S6").Value = 0.00005311868906042436
$ws.Range("T6").Value = 0.00005615115775998398
$ws.Range("I7").Value = 0.9827953701592058
$ws.Range("J7").Value = 0.9827953701592059
$ws.Range("M7").Value = 1.824475333333333
$ws.Range("N7").Value = 5.473426
$ws.Range("O7").Value = 0.1906606574278047
$ws.Range("P7").Value = 0.2015451970524477
$ws.Range("Q7").Value = 9.039167820139332
$ws.Range("R7").Value = 81.35251038125399
$ws.Range("S7").Value = 0.1873804113915569
$ws.Range("T7").Value = 0.1980776865409704
$ws.Range("I8").Value = 0.9827953701592058
$ws.Range("J8").Value = 0.9827953701592059
$ws.Range("O8").Value = 0.6423822165107047
$ws.Range("P8").Value = 0.6790548829333741
$ws.Range("S8").Value = 0.6313302682593291
$ws.Range("T8").Value = 0.6673719950309217
$ws.Range("I9").Value = 0.9827953701592058
$ws.Range("J9").Value = 0.9827953701592059
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.009389666666666666
$ws.Range("N9").Value = 0.028169
$ws.Range("O9").Value = 0.0009812355294625031
$ws.Range("P9").Value = 0.001037252838673693
$ws.Range("Q9").Value = 0.04652009880566666
$ws.Range("R9").Value = 0.4186808892509999
$ws.Range("S9").Value = 0.0009643537353914651
$ws.Range("T9").Value = 0.001019407287533
$ws.Range("D10").Value = "MuSCs"
$ws.Range("I10").Value = 0.9827953701592058
$ws.Range("J10").Value = 0.9827953701592059
$ws.Range("M10").Value = 1.5503715
$ws.Range("N10").Value = 3.100743
$ws.Range("O10").Value = 0.1620163583726162
$ws.Range("P10").Value = 0.1141770910840848
$ws.Range("Q10").Value = 7.681149706999499
$ws.Range("R10").Value = 46.086898241997
$ws.Range("S10").Value = 0.1592289268986619
$ws.Range("T10").Value = 0.1122127164956845
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("I11").Value = 0.9827953701592058
$ws.Range("J11").Value = 0.9827953701592059
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.03788966666666667
$ws.Range("N11").Value = 0.113669
$ws.Range("O11").Value = 0.003959532159411881
$ws.Range("P11").Value = 0.004185576091419648
$ws.Range("Q11").Value = 0.1877202993056667
$ws.Range("R11").Value = 1.689482693751
$ws.Range("S11").Value = 0.003891409874266479
$ws.Range("T11").Value = 0.004113564804096295
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.019102
$ws.Range("H12").Value = 0.057306
$ws.Range("I12").Value = 0.003789234556237495
$ws.Range("J12").Value = 0.003789234556237496
$ws.Range("M12").Value = 1.824475333333333
$ws.Range("N12").Value = 5.473426
$ws.Range("O12").Value = 0.1906606574278047
$ws.Range("P12").Value = 0.2015451970524477
$ws.Range("Q12").Value = 0.03485112781733334
$ws.Range("R12").Value = 0.313660150356
$ws.Range("S12").Value = 0.0007224579516403967
$ws.Range("T12").Value = 0.0007637020253148302
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.019102
$ws.Range("H13").Value = 0.057306
$ws.Range("I13").Value = 0.003789234556237495
$ws.Range("J13").Value = 0.003789234556237496
$ws.Range("O13").Value = 0.6423822165107047
$ws.Range("P13").Value = 0.6790548829333741
$ws.Range("Q13").Value = 0.117421942404
$ws.Range("R13").Value = 1.056797481636
$ws.Range("S13").Value = 0.002434136893114799
$ws.Range("T13").Value = 0.002573098227992949
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.019102
$ws.Range("H14").Value = 0.057306
$ws.Range("I14").Value = 0.003789234556237495
$ws.Range("J14").Value = 0.003789234556237496
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.009389666666666666
$ws.Range("N14").Value = 0.028169
$ws.Range("O14").Value = 0.0009812355294625031
$ws.Range("P14").Value = 0.001037252838673693
$ws.Range("Q14").Value = 0.0001793614126666667
$ws.Range("R14").Value = 0.001614252714
$ws.Range("S14").Value = 0.000003718131576047312
$ws.Range("T14").Value = 0.000003930394299857796
$ws.Range("D15").Value = "MuSCs"
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.019102
$ws.Range("H15").Value = 0.057306
$ws.Range("I15").Value = 0.003789234556237495
$ws.Range("J15").Value = 0.003789234556237496
$ws.Range("M15").Value = 1.5503715
$ws.Range("N15").Value = 3.100743
$ws.Range("O15").Value = 0.1620163583726162
$ws.Range("P15").Value = 0.1141770910840848
$ws.Range("Q15").Value = 0.029615196393
$ws.Range("R15").Value = 0.177691178358
$ws.Range("S15").Value = 0.0006139179838212754
$ws.Range("T15").Value = 0.0004326437790664902
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.019102
$ws.Range("H16").Value = 0.057306
$ws.Range("I16").Value = 0.003789234556237495
$ws.Range("J16").Value = 0.003789234556237496
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.03788966666666667
$ws.Range("N16").Value = 0.113669
$ws.Range("O16").Value = 0.003959532159411881
$ws.Range("P16").Value = 0.004185576091419648
$ws.Range("Q16").Value = 0.0007237684126666667
$ws.Range("R16").Value = 0.006513915714000001
$ws.Range("S16").Value = 0.00001500359608497717
$ws.Range("T16").Value = 0.0000158601295633688

Write-Output "Updated cells successfully"
